$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 362.26
$ws.Range("I15").Value = 362.26
$ws.Range("K15").Value = 1086.78
$ws.Range("M15").Value = -917.78
$ws.Range("H76").Value = 3262.8918
$ws.Range("I76").Value = 2889.6
$ws.Range("K76").Value = 2889.6
$ws.Range("M76").Value = -2574.6
$ws.Range("H79").Value = 3262.8918
$ws.Range("I79").Value = 2889.6
$ws.Range("K79").Value = 2889.6
$ws.Range("M79").Value = -1797.6
$ws.Range("H137").Value = 13891696
$ws.Range("I137").Value = 25001652
$ws.Range("J137").Value = 4250.375
$ws.Range("K137").Value = 75004956
$ws.Range("L137").Value = 12751.125
$ws.Range("M137").Value = -75002406
$ws.Range("N137").Value = -17851.125
$ws.Range("H138").Value = 2223.6628
$ws.Range("I138").Value = 1045.9286
$ws.Range("J138").Value = 2452.6667
$ws.Range("K138").Value = 3137.7858
$ws.Range("L138").Value = 7358.000100000001
$ws.Range("M138").Value = 2002.2142
$ws.Range("N138").Value = -17638.0001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 12277.556
$ws.Range("J54").Value = 12312.25
$ws.Range("L54").Value = 12312.25
$ws.Range("N54").Value = -13850.25
$ws.Range("H74").Value = 1720.279
$ws.Range("I74").Value = 1422
$ws.Range("J74").Value = 4628.5
$ws.Range("K74").Value = 1422
$ws.Range("L74").Value = 4628.5
$ws.Range("M74").Value = -548
$ws.Range("N74").Value = -6376.5
$ws.Range("H77").Value = 1720.279
$ws.Range("I77").Value = 1422
$ws.Range("J77").Value = 4628.5
$ws.Range("K77").Value = 7110
$ws.Range("L77").Value = 23142.5
$ws.Range("M77").Value = -2742
$ws.Range("N77").Value = -31878.5
$ws.Range("H104").Value = 17681.25
$ws.Range("J104").Value = 17681.25
$ws.Range("L104").Value = 17681.25
$ws.Range("N104").Value = -24669.25
$ws.Range("H132").Value = 1279.2128
$ws.Range("I132").Value = 699.94446
$ws.Range("K132").Value = 2099.83338
$ws.Range("M132").Value = 430.16662
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9920
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 9920
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 9920
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -10144
$ws.Range("H28").Value = 14982.25
$ws.Range("J28").Value = 14982.25
$ws.Range("L28").Value = 14982.25
$ws.Range("N28").Value = -15472.25
$ws.Range("H132").Value = 1737.4412
$ws.Range("I132").Value = 1366.44
$ws.Range("J132").Value = 2768
$ws.Range("K132").Value = 4099.32
$ws.Range("L132").Value = 8304
$ws.Range("M132").Value = -1569.32
$ws.Range("N132").Value = -13364
$ws.Range("H134").Value = 2913.375
$ws.Range("I134").Value = 1722.95
$ws.Range("J134").Value = 4897.4165
$ws.Range("K134").Value = 5168.85
$ws.Range("L134").Value = 14692.2495
$ws.Range("M134").Value = -2633.85
$ws.Range("N134").Value = -19762.2495
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 476495.2
$ws.Range("I107").Value = 166.83333
$ws.Range("J107").Value = 667026.5600000001
$ws.Range("K107").Value = 500.49999
$ws.Range("L107").Value = 2001079.68
$ws.Range("M107").Value = 1419.50001
$ws.Range("N107").Value = -2004919.68
$ws.Range("H114").Value = 3071.625
$ws.Range("I114").Value = 2746.5
$ws.Range("J114").Value = 3180
$ws.Range("K114").Value = 8239.5
$ws.Range("L114").Value = 9540
$ws.Range("M114").Value = -4985.5
$ws.Range("N114").Value = -16048
$ws.Range("H121").Value = 5883127
$ws.Range("I121").Value = 403.75
$ws.Range("J121").Value = 11112215
$ws.Range("K121").Value = 1211.25
$ws.Range("L121").Value = 33336645
$ws.Range("M121").Value = 98.75
$ws.Range("N121").Value = -33339265
$ws.Range("H129").Value = 4754.1714
$ws.Range("I129").Value = 1792.8462
$ws.Range("J129").Value = 6504.0454
$ws.Range("K129").Value = 5378.5386
$ws.Range("L129").Value = 19512.1362
$ws.Range("M129").Value = -378.5385999999999
$ws.Range("N129").Value = -29512.1362
$ws.Range("H131").Value = 843.15
$ws.Range("I131").Value = 382
$ws.Range("J131").Value = 909.02856
$ws.Range("K131").Value = 1146
$ws.Range("L131").Value = 2727.08568
$ws.Range("M131").Value = 3894
$ws.Range("N131").Value = -12807.08568
$ws.Range("H132").Value = 6585689
$ws.Range("I132").Value = 1200
$ws.Range("J132").Value = 7408750
$ws.Range("K132").Value = 10800
$ws.Range("L132").Value = 66678750
$ws.Range("M132").Value = -8270
$ws.Range("N132").Value = -66683810
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 5175
$ws.Range("I5").Value = 350
$ws.Range("J5").Value = 10000
$ws.Range("K5").Value = 350
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = -238
$ws.Range("N5").Value = -10224
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9310.714
$ws.Range("J2").Value = 9310.714
$ws.Range("L2").Value = 9310.714
$ws.Range("N2").Value = -9534.714
$ws.Range("H104").Value = 20274
$ws.Range("J104").Value = 20274
$ws.Range("L104").Value = 20274
$ws.Range("N104").Value = -27262
$ws.Range("H122").Value = 2025.5
$ws.Range("I122").Value = 2067.3333
$ws.Range("K122").Value = 6201.999899999999
$ws.Range("M122").Value = -3751.999899999999
$ws.Range("H128").Value = 36085.6
$ws.Range("J128").Value = 36085.6
$ws.Range("L128").Value = 36085.6
$ws.Range("N128").Value = -46045.6
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 6283797.5
$ws.Range("J2").Value = 6283797.5
$ws.Range("L2").Value = 6283797.5
$ws.Range("N2").Value = -6284021.5
$ws.Range("H10").Value = 48365.6
$ws.Range("I10").Value = 905
$ws.Range("J10").Value = 80006
$ws.Range("K10").Value = 905
$ws.Range("L10").Value = 80006
$ws.Range("M10").Value = -736
$ws.Range("N10").Value = -80344
$ws.Range("H51").Value = 3500
$ws.Range("I51").Value = 3500
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 3500
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -2990
$ws.Range("N51").ClearContents()
$ws.Range("H132").Value = 2554.5386
$ws.Range("I132").Value = 1974.5264
$ws.Range("J132").Value = 3105.55
$ws.Range("K132").Value = 5923.5792
$ws.Range("L132").Value = 9316.650000000001
$ws.Range("M132").Value = -3393.5792
$ws.Range("N132").Value = -14376.65
